# Update countries & provincias Spain
# - Refreshed COVID-19 case counts for several countries (Alemania, Taiwan,
#   and the rows that move because Bulgaria / El Salvador overtook their
#   neighbours in the total-cases ranking), and refreshed the "last updated"
#   timestamp.
# - Because the list is kept sorted by total cases, Bulgaria's new total
#   (1488) now outranks Cuba (1467) and Nueva Zelanda (1476), so those two
#   rows shift down one slot; El Salvador's new total (395) now outranks
#   Kenia (384), so they swap places too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp footer (row 1)
$ws.Range("A1").Value = "Datos actualizados a 30 de Abril de 2020 a las 08:22"

# Alemania (row 9) - refreshed totals
$ws.Range("D9").Value = 123500
$ws.Range("E9").Value = 31572

# Rows 79-81: Bulgaria jumps ahead of Nueva Zelanda and Cuba
$ws.Range("A79").Value = "Bulgaria"
$ws.Range("B79").Value = 1488
$ws.Range("C79").Value = 41
$ws.Range("D79").Value = 266
$ws.Range("E79").Value = 1157
$ws.Range("F79").Value = 38
$ws.Range("G79").Value = 1
$ws.Range("H79").Value = 65

$ws.Range("A80").Value = "Nueva Zelanda"
$ws.Range("B80").Value = 1476
$ws.Range("C80").Value = 2
$ws.Range("D80").Value = 1241
$ws.Range("E80").Value = 216
$ws.Range("F80").Value = 1
$ws.Range("G80").Value = 0
$ws.Range("H80").Value = 19

$ws.Range("A81").Value = "Cuba"
$ws.Range("B81").Value = 1467
$ws.Range("C81").Value = 0
$ws.Range("D81").Value = 617
$ws.Range("E81").Value = 792
$ws.Range("F81").Value = 14
$ws.Range("G81").Value = 0
$ws.Range("H81").Value = 58

# Taiwan (row 116) - refreshed totals
$ws.Range("D116").Value = 322
$ws.Range("E116").Value = 101

# Rows 119-120: El Salvador jumps ahead of Kenia
$ws.Range("A119").Value = "El Salvador"
$ws.Range("B119").Value = 395
$ws.Range("C119").Value = 18
$ws.Range("D119").Value = 118
$ws.Range("E119").Value = 268
$ws.Range("F119").Value = 3
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 9

$ws.Range("A120").Value = "Kenia"
$ws.Range("B120").Value = 384
$ws.Range("C120").Value = 0
$ws.Range("D120").Value = 129
$ws.Range("E120").Value = 240
$ws.Range("F120").Value = 2
$ws.Range("G120").Value = 0
$ws.Range("H120").Value = 15
